$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "Approximately $5.5 per packet" - runs already contain the right text, just
#    force a merge of the split runs via a no-op replace.
Replace-Text "Approximately $5.5 per packet" "Approximately $5.5 per packet"

# 2. "Price: $33" - merge split runs (no textual change).
Replace-Text "Price: $33" "Price: $33"

# 3. "Red Onion 2 pecies" -> add "(one big and one small)"
Replace-Text "Red Onion 2 pecies" "Red Onion 2 pecies (one big and one small)"

# 4. "Approximately $4 for 1" -> "Approximately $7 for 2"
Replace-Text "Approximately $4 for 1" "Approximately $7 for 2"

# 5. "Price: $8" -> "Price: $7"
Replace-Text "Price: $8" "Price: $7"

# 6. "Price: $29" - merge split runs (no textual change).
Replace-Text "Price: $29" "Price: $29"

# 7. "Pigeon Brand Fermented Vegetable 140g 3 Can" -> "Vegetable 2"
Replace-Text "Pigeon Brand Fermented Vegetable 140g 3 Can" "Vegetable 2"

# 8. "Approximately $7.9 per Can" -> "Approximately $11.9 per Can"
Replace-Text "Approximately $7.9 per Can" "Approximately $11.9 per Can"

# 9. "Price: $23.7" -> "Price: $23.8"
Replace-Text "Price: $23.7" "Price: $23.8"

# 10. Total price "300" -> "299.1"
Replace-Text "Total Price: Approximate $300 (All in terms of Hong Kong Dollars)" "Total Price: Approximate $299.1 (All in terms of Hong Kong Dollars)"

# 11. "Updated: 19" - merge split runs (no textual change).
Replace-Text "Updated: 19" "Updated: 19"

# 12. "Beijing Noodle with Fermented Vegetable" -> "Beijing Noodle with Vegetable"
Replace-Text "Beijing Noodle with Fermented Vegetable –  (30 meals) (Served for Breakfast or Lunch)" "Beijing Noodle with Vegetable –  (30 meals) (Served for Breakfast or Lunch)"

# 13. "Add Chicken Powder." - merge split runs (no textual change).
Replace-Text "Add Chicken Powder." "Add Chicken Powder."

# 14. "1/2 Onion." - merge split runs (no textual change).
Replace-Text "1/2 Onion." "1/2 Onion."
